$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TableTK_Seq_Metrics")

# Rename the "File Size" header (column E) to be explicit about the file type.
$ws.Range("E3").Value = "File Size (fasta)"

# Fill in newly-gathered raw file stats for the 10X Illumina runs (rows 9-12).

# Row 9: Female, 10X Illumina, R1
$ws.Range("I9").Value = 151
$ws.Range("J9").Value = 33

# Row 10: Female, 10X Illumina, R2
$ws.Range("I10").Value = 151
$ws.Range("J10").Value = 32

# Row 11: Male, 10X Illumina, R1
$ws.Range("F11").Value = 65806680934
$ws.Range("F11").NumberFormat = "#,##0"
$ws.Range("H11").Value = 435805834
$ws.Range("H11").NumberFormat = "#,##0"
$ws.Range("I11").Value = 151
$ws.Range("J11").Value = 34

# Row 12: Male, 10X Illumina, R2
$ws.Range("F12").Value = 65806680934
$ws.Range("F12").NumberFormat = "#,##0"
$ws.Range("H12").Value = 435805834
$ws.Range("H12").NumberFormat = "#,##0"
$ws.Range("I12").Value = 151
$ws.Range("J12").Value = 32

# Leave the cursor where the author last left it on this sheet.
$ws.Activate() | Out-Null
$ws.Range("F20").Select() | Out-Null
